# Update "Total Users" and "Administrators" counts in the summary tables.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# First summary table (rows 3-4): A4 = Total Users, B4 = Administrators
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = 2

# "Total Statistics" table (rows 17-18): C18 = Total Users, D18 = Administrators
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 2
